$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date rolled forward from 06-16 to 06-17)
$ws.Name = "Through 2022-06-17"

# Update the "June (through 06-16)" label to "June (through 06-17)"
$ws.Range("A7").Value = "June (through 06-17)"

# Row 7 - June counts by year (B..I = 2015..2022)
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 34
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = 29
$ws.Range("G7").Value = 63
$ws.Range("H7").Value = 63
$ws.Range("I7").Value = 82

# Row 8 - Total counts by year (B..I = 2015..2022)
$ws.Range("B8").Value = 117
$ws.Range("C8").Value = 230
$ws.Range("D8").Value = 350
$ws.Range("E8").Value = 330
$ws.Range("F8").Value = 233
$ws.Range("G8").Value = 421
$ws.Range("H8").Value = 694
$ws.Range("I8").Value = 745
